$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.25955867767334
$ws.Range("B1").Value = 1.684051871299744
$ws.Range("C1").Value = 5.542348384857178
$ws.Range("D1").Value = 1.556192874908447
$ws.Range("E1").Value = 0.7474992871284485
